# Automated data refresh: updates extraction timestamps and latest
# meteorological readings (humidity, precipitation, temperature,
# pressure, snow depth, wind) for the rows whose source values changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-09 19:48:33"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "91%"
$ws.Range("E3").Value = "2026-02-09 19:48:35"
$ws.Range("I3").Value = "1.5 mm"
$ws.Range("E4").Value = "2026-02-09 19:48:38"
$ws.Range("O4").Value = "8.3 °C"
$ws.Range("E5").Value = "2026-02-09 19:48:40"
$ws.Range("E6").Value = "2026-02-09 19:48:43"
$ws.Range("E7").Value = "2026-02-09 19:48:45"
$ws.Range("E8").Value = "2026-02-09 19:48:48"
$ws.Range("E9").Value = "2026-02-09 19:48:51"
$ws.Range("E10").Value = "2026-02-09 19:48:53"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "83%"
$ws.Range("O10").Value = "8.2 °C"
$ws.Range("E11").Value = "2026-02-09 19:48:56"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "82%"
$ws.Range("E12").Value = "2026-02-09 19:48:58"
$ws.Range("E13").Value = "2026-02-09 19:49:01"
$ws.Range("E14").Value = "2026-02-09 19:49:04"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "76%"
$ws.Range("E15").Value = "2026-02-09 19:49:06"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "80%"
$ws.Range("O15").Value = "8.3 °C"
$ws.Range("E16").Value = "2026-02-09 19:49:09"
$ws.Range("E17").Value = "2026-02-09 19:49:12"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "82%"
$ws.Range("E18").Value = "2026-02-09 19:49:14"
$ws.Range("J18").Value = "1007.2 hPa"
$ws.Range("O18").Value = "9.1 °C"
$ws.Range("E19").Value = "2026-02-09 19:49:17"
$ws.Range("E20").Value = "2026-02-09 19:49:20"
$ws.Range("O20").Value = "-4.4 °C"
$ws.Range("E21").Value = "2026-02-09 19:49:22"
$ws.Range("E22").Value = "2026-02-09 19:49:25"
$ws.Range("E23").Value = "2026-02-09 19:49:28"
$ws.Range("O23").Value = "-3.6 °C"
$ws.Range("E24").Value = "2026-02-09 19:49:30"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "82%"
$ws.Range("I24").Value = "0.7 mm"
$ws.Range("E25").Value = "2026-02-09 19:49:33"
$ws.Range("E26").Value = "2026-02-09 19:49:35"
$ws.Range("G26").Value = "2 cm"
$ws.Range("E27").Value = "2026-02-09 19:49:37"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "82%"
$ws.Range("L27").Value = "23.0 km/h - 265º 19:10 TU"
$ws.Range("O27").Value = "-2.4 °C"
$ws.Range("E28").Value = "2026-02-09 19:49:40"
$ws.Range("E29").Value = "2026-02-09 19:49:42"
$ws.Range("E30").Value = "2026-02-09 19:49:45"
$ws.Range("E31").Value = "2026-02-09 19:49:47"
$ws.Range("E32").Value = "2026-02-09 19:49:49"
$ws.Range("E33").Value = "2026-02-09 19:49:52"
$ws.Range("O33").Value = "3.0 °C"
$ws.Range("E34").Value = "2026-02-09 19:49:55"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "75%"
$ws.Range("E35").Value = "2026-02-09 19:49:57"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "77%"
$ws.Range("I35").Value = "1.4 mm"
$ws.Range("E36").Value = "2026-02-09 19:50:00"
$ws.Range("E37").Value = "2026-02-09 19:50:03"
$ws.Range("E38").Value = "2026-02-09 19:50:06"
$ws.Range("E39").Value = "2026-02-09 19:50:08"
$ws.Range("O39").Value = "-3.4 °C"
$ws.Range("E40").Value = "2026-02-09 19:50:10"
$ws.Range("E41").Value = "2026-02-09 19:50:13"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "56%"
$ws.Range("E42").Value = "2026-02-09 19:50:15"
$ws.Range("E43").Value = "2026-02-09 19:50:18"
$ws.Range("E44").Value = "2026-02-09 19:50:21"
$ws.Range("E45").Value = "2026-02-09 19:50:23"
$ws.Range("J45").Value = "1007.2 hPa"
$ws.Range("E46").Value = "2026-02-09 19:50:26"
$ws.Range("I46").Value = "0.6 mm"
$ws.Range("J46").Value = "1008.8 hPa"
